# Documentation updates: fixed links, typos, versions, etc. (#410)
#
# The "Inventory" worksheet documents the repo's directory/file layout as a
# markdown table. This change:
#   - Renames the "Gradle Properties" doc link to "Tenant Configuration"
#     (G8), and removes the now-redundant "JavaScript Template Files" doc
#     link cell (G22), since that callout has been folded into F22's text.
#   - Updates a couple of descriptions to reflect the base ml-config
#     directory applying to "tenants" (not "environments"), and to expand
#     on what the local, non-admin deployment user is used for.
#   - Rewords the /templates description to note it's reserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")

$ws.Range("G8").Value = "[Tenant Configuration](/docs/lux-backend-deployment.md#tenant-configuration)"

$ws.Range("F9").Value = "The base configuration directory applicable to all tenants.  It includes the group configuration, main content database, roles, and application servers."

$ws.Range("F11").Value = "Intended for local developer environments.  Defines a local, non-admin user to perform most of deployments with plus some endpoint consumers."

$ws.Range("F22").Value = "Reserved for JavaScript template files used by [/build.gradle](/build.gradle)."
$ws.Range("G22").ClearContents()
